# "Serach on 16 Sept" - refresh the local job search results on Page1
# with the latest scrape: some previous listings fell off, new ones
# appear, and the remaining ones shift position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page1")

$rows = @(
  @{ Title = "English Turkish Language Expert"; Url = "https://tr.indeed.com/rc/clk?jk=d1c793aff89a89c3&fccid=b31babcb6321f38f&vjs=3" },
  @{ Title = "Senior Software QA Engineer"; Url = "https://tr.indeed.com/rc/clk?jk=3f833354847da52e&fccid=83cbd453678ade10&vjs=3" },
  @{ Title = "QA Engineer, IgniteTech (Remote) - `$30,000/year USD"; Url = "https://tr.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0Bp-pWyh_yd7-BMq_HjiKF8q-wDEEHI6K3EhZcTze0zEvR-NFCW1RxZBIlN3zMxg-_a0rRewiNcEPttf6EywkUo6jfD0QBzPdWMed2TkhoCyMfZxBILCVZY5YccN9_KvR2JMS0lFuAD8BsBDLICcUmxLP6Yv3Z0nxVsNLNy0dO8_K77dJjpASCrp-ruwTmPx8PJ9z0Vl9ACvO78lECpkVLpn9wU9e_2AuMKWP70c-P3OKWym0HZvQgofgV8qnjXNk2FEv1Lqr9QZld4ZI7s0yWZV343pfFeL3obBhw8hDL1T0EcOX_j5kG1YHsj2Ec7xDaTibUBo4y2Xt2QfUgnP4oJ6zTdLIoMqikcn3vQbdxNu7AqexQbJumVrB-WuFdrzhNPiBGjNi3LXZIoyEizCRP37ibpEZ2ghyBBsWmjeIJGPtm3WQRPX8Ou_cbtN1uNzKGnUXdSJoMVLf2OwqWthlATjoQXvAasNJoPPeUdoHVOMZNqS3zpbG7MhXgHhemLh6gqn9K3mF4Qpex8gGRBZ0VXgvb9VBJiD2FUIA4I7x6bIFdl6o5TGDT5Eg71M_O42be--um2R8PdhtVTWlW7rgKS-agODUak_nlwRi1WqaipaTFf1DVPKgCO&xkcb=SoD_-_M3ZjL9eiR9KZ0JbzkdCdPP&p=2&fvj=0&vjs=3" },
  @{ Title = "QA Manual Engineer"; Url = "https://tr.indeed.com/rc/clk?jk=b6ee03f1b74fa59a&fccid=da69cc7fb1190c2f&vjs=3" },
  @{ Title = "UI Developer"; Url = "https://tr.indeed.com/rc/clk?jk=87cb4e48c95b2206&fccid=a1247e0f0273002a&vjs=3" },
  @{ Title = "Senior Software Test Engineer - Robotics"; Url = "https://tr.indeed.com/rc/clk?jk=f5cd9bbb6644a4f8&fccid=bc7b1cb755d12a77&vjs=3" },
  @{ Title = "QA Storage Test Expert"; Url = "https://tr.indeed.com/rc/clk?jk=e42205987b070b3d&fccid=404d8dccef12ec23&vjs=3" },
  @{ Title = "Dedicated Engineers"; Url = "https://tr.indeed.com/rc/clk?jk=a2fc83080e802ff0&fccid=db5bd09bf2cdc268&vjs=3" },
  @{ Title = "Project Quality Leader"; Url = "https://tr.indeed.com/rc/clk?jk=026c9514f10b9c41&fccid=1a58d1f68a45eff2&vjs=3" },
  @{ Title = "Senior Software Engineer"; Url = "https://tr.indeed.com/rc/clk?jk=75ffe91be437383d&fccid=389b8a90ef238b4f&vjs=3" },
  @{ Title = "Lead Supplier Quality Inspector (Mechanical Commodity)"; Url = "https://tr.indeed.com/rc/clk?jk=feaf2936129955a7&fccid=1a58d1f68a45eff2&vjs=3" },
  @{ Title = "In House CRA"; Url = "https://tr.indeed.com/rc/clk?jk=95853af20ebefb35&fccid=6c527578251bc5e1&vjs=3" },
  @{ Title = "Director of Chaos Engineering - Ubuntu & Distributed Systems"; Url = "https://tr.indeed.com/rc/clk?jk=473c4690d08d8592&fccid=3911cf63fc0692f4&vjs=3" }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $rows[$i].Title
  $ws.Cells.Item($r, 2).Value = $rows[$i].Url
}
